# InstantiateAgents.xlsx - "Model Test with larger init CSVs"
# Adds a second data row to each of the three init sheets (FixNodes,
# VarNodes, Transporters) plus a third row to Transporters, and moves
# the active sheet/selection from Transporters to FixNodes.

$wb = $excel.ActiveWorkbook

$wsFix   = $wb.Worksheets.Item("FixNodes")
$wsVar   = $wb.Worksheets.Item("VarNodes")
$wsTrans = $wb.Worksheets.Item("Transporters")

# --- Transporters: two new rows (sheet rows 3 and 4) ---
# (filled in first so the new shared strings are interned in the same
# order as the original authored workbook: TestTrans2, VarTest2,
# FixTest2, RedOrigin, TestTrans3, ChadGalactic)
$wsTrans.Range("A3").Value = "TestTrans2"
$wsTrans.Range("B3").Value = 1000
$wsTrans.Range("C3").Value = 350
$wsTrans.Range("D3").Value = 120
$wsTrans.Range("E3").Value = "VarTest2"
$wsTrans.Range("F3").Value = "FixTest2"
$wsTrans.Range("G3").Value = "RedOrigin"

$wsTrans.Range("A4").Value = "TestTrans3"
$wsTrans.Range("B4").Value = 0
$wsTrans.Range("C4").Value = 300
$wsTrans.Range("D4").Value = 200
$wsTrans.Range("E4").Value = "earth"
$wsTrans.Range("F4").Value = "VarTest1"
$wsTrans.Range("G4").Value = "ChadGalactic"

# --- FixNodes: new row 2 (sheet row 3) ---
$wsFix.Range("A3").Value = "FixTest2"
$wsFix.Range("B3").Value = 10000
$wsFix.Range("C3").Value = 25000
$wsFix.Range("D3").Value = 4000
$wsFix.Range("E3").Value = 12

# --- VarNodes: new row 2 (sheet row 3) ---
$wsVar.Range("A3").Value = "VarTest2"
$wsVar.Range("B3").Value = 800
$wsVar.Range("C3").Value = 950
$wsVar.Range("D3").Value = 50
$wsVar.Range("E3").Value = 2

# Column G on Transporters now needs a best-fit width for the longer
# strings that were just added (id column G, e.g. "ChadGalactic").
$wsTrans.Columns.Item(7).ColumnWidth = 10.830729166666666

# Update selections on the non-active sheets first ...
[void]$wsVar.Range("D21").Select()
[void]$wsTrans.Range("A5").Select()

# ... then move the active sheet/selection to FixNodes, matching the
# workbook's new activeTab / tabSelected state.
[void]$wsFix.Select()
[void]$wsFix.Range("D4").Select()
